$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header region values
$ws.Range("B1").Value = 4
$ws.Range("B2").Value = $null

# Add new student rows
$ws.Range("A4").Value = 44
$ws.Range("B4").Value = "SAIDI"
$ws.Range("C4").Value = "SEIF"
$ws.Range("D4").Value = "NGONGO"
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = 45
$ws.Range("B5").Value = "SALMA"
$ws.Range("C5").Value = "SEIF"
$ws.Range("D5").Value = "NGONGO"
$ws.Range("E5").Value = 2
